$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.72"
$ws.Range("E2").Value = "'-4.87%"
$ws.Range("D3").Value = "'31.61"
$ws.Range("E3").Value = "'-1.13%"
$ws.Range("D4").Value = "'5.132"
$ws.Range("E4").Value = "'-3.95%"
$ws.Range("D5").Value = "'0.07504"
$ws.Range("E5").Value = "'-0.89%"
$ws.Range("D6").Value = "'7.738"
$ws.Range("E6").Value = "'-1.34%"
$ws.Range("E7").Value = "'5.87%"
$ws.Range("D8").Value = "'3.799"
$ws.Range("E8").Value = "'2.43%"
$ws.Range("D9").Value = "'0.9319"
$ws.Range("E9").Value = "'1.74%"
$ws.Range("D10").Value = "'0.1694"
$ws.Range("E10").Value = "'-1.46%"
$ws.Range("D11").Value = "'0.07163"
$ws.Range("E11").Value = "'-6.79%"
$ws.Range("D12").Value = "'0.07921"
$ws.Range("E12").Value = "'-3.88%"
$ws.Range("E13").Value = "'-0.72%"
$ws.Range("D14").Value = "'0.09904"
$ws.Range("E14").Value = "'0.22%"
$ws.Range("D15").Value = "'0.001501"
$ws.Range("E15").Value = "'-1.95%"
$ws.Range("D16").Value = "'0.006372"
$ws.Range("E16").Value = "'-0.80%"
$ws.Range("D17").Value = "'3.445"
$ws.Range("E17").Value = "'-1.26%"
$ws.Range("D18").Value = "'2.226"
$ws.Range("E18").Value = "'-0.58%"
$ws.Range("D19").Value = "'0.3284"
$ws.Range("E19").Value = "'-0.89%"
$ws.Range("E20").Value = "'0.84%"
$ws.Range("D21").Value = "'4.577"
$ws.Range("E21").Value = "'9.47%"
$ws.Range("D22").Value = "'0.04664"
$ws.Range("E22").Value = "'2.35%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-0.16%"
$ws.Range("D25").Value = "'0.004422"
$ws.Range("E25").Value = "'-1.68%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.17%"
$ws.Range("D27").Value = "'0.0001878"
$ws.Range("E27").Value = "'7.88%"
$ws.Range("D39").Value = "'0.01673"
$ws.Range("E39").Value = "'-1.26%"
$ws.Range("D40").Value = "'0.04456"
$ws.Range("E40").Value = "'-3.22%"
$ws.Range("D41").Value = "'0.007065"
$ws.Range("E41").Value = "'-2.58%"
$ws.Range("E42").Value = "'-3.02%"
$ws.Range("D43").Value = "'0.002063"
$ws.Range("E43").Value = "'-8.70%"
$ws.Range("D44").Value = "'0.01128"
$ws.Range("E44").Value = "'-20.06%"
$ws.Range("D45").Value = "'0.00006017"
$ws.Range("E45").Value = "'-2.36%"
$ws.Range("D46").Value = "'1.918"
$ws.Range("E46").Value = "'1.34%"
$ws.Range("E47").Value = "'-0.24%"
